# update get otp api
$wb = $excel.ActiveWorkbook

$wsForgot = $wb.Worksheets.Item("ForgotPinOtp")
$wsGetOtp = $wb.Worksheets.Item("GetOtp")

# Old row 6 (empty cell, A6 s="3") is removed; the old row 7 (A7 s="6", value 181)
# shifts up into row 6, carrying its style along -> used range shrinks to A1:A6.
$wsGetOtp.Rows(6).Delete()

# The shifted-up cell's value changes from 181 to 182.
$wsGetOtp.Range("A6").Value = 182

# GetOtp becomes the active sheet/tab, with its selection moved to D11.
$wsGetOtp.Activate()
$wsGetOtp.Range("D11").Select()

# ForgotPinOtp's selection stays at D24, but it's no longer the active tab.
$wsForgot.Range("D24").Select()
$wsGetOtp.Activate()
